$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 12463
$ws.Range("E2").Value = 1855
$ws.Range("F2").Value = 1855
$ws.Range("G2").Value = 1847
$ws.Range("H2").Value = 1433
$ws.Range("I2").Value = 1281
$ws.Range("J2").Value = 152
$ws.Range("K2").Value = 13941
$ws.Range("L2").Value = 3818
$ws.Range("M2").Value = 10123
$ws.Range("N2").Value = 9722
$ws.Range("O2").Value = 401
$ws.Range("P2").Value = 222
$ws.Range("Q2").Value = 1947
$ws.Range("R2").Value = -929
$ws.Range("S2").Value = -680
$ws.Range("T2").Value = 1173
$ws.Range("U2").Value = 774
$ws.Range("V2").Value = 1106
$ws.Range("W2").Value = 14.88
$ws.Range("X2").Value = 11.5
$ws.Range("Y2").Value = 14.05
$ws.Range("Z2").Value = 11
$ws.Range("AA2").Value = 37.72
$ws.Range("AB2").Value = 4473.05
$ws.Range("AC2").Value = 2891
$ws.Range("AD2").Value = 18.33
$ws.Range("AE2").Value = 21940
$ws.Range("AF2").Value = 2.42
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 0.38
$ws.Range("AI2").Value = 6.92
$ws.Range("AJ2").Value = 44311468

# Row 3
$ws.Range("D3").Value = 15849
$ws.Range("E3").Value = 1968
$ws.Range("F3").Value = 1968
$ws.Range("G3").Value = 1908
$ws.Range("H3").Value = 1445
$ws.Range("I3").Value = 1299
$ws.Range("J3").Value = 146
$ws.Range("K3").Value = 22210
$ws.Range("L3").Value = 8586
$ws.Range("M3").Value = 13624
$ws.Range("N3").Value = 11767
$ws.Range("O3").Value = 1858
$ws.Range("P3").Value = 222
$ws.Range("Q3").Value = 2584
$ws.Range("R3").Value = -1671
$ws.Range("S3").Value = 1564
$ws.Range("T3").Value = 564
$ws.Range("U3").Value = 2020
$ws.Range("V3").Value = 3800
$ws.Range("W3").Value = 12.42
$ws.Range("X3").Value = 9.119999999999999
$ws.Range("Y3").Value = 12.09
$ws.Range("Z3").Value = 7.99
$ws.Range("AA3").Value = 63.02
$ws.Range("AB3").Value = 5019.32
$ws.Range("AC3").Value = 2931
$ws.Range("AD3").Value = 14.36
$ws.Range("AE3").Value = 26554
$ws.Range("AF3").Value = 1.59
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 0.48
$ws.Range("AI3").Value = 6.82
$ws.Range("AJ3").Value = 44311468

# Row 4
$ws.Range("D4").Value = 20016
$ws.Range("E4").Value = 1794
$ws.Range("F4").Value = 1794
$ws.Range("G4").Value = 1538
$ws.Range("H4").Value = 1090
$ws.Range("I4").Value = 1100
$ws.Range("J4").Value = -10
$ws.Range("K4").Value = 23897
$ws.Range("L4").Value = 10122
$ws.Range("M4").Value = 13775
$ws.Range("N4").Value = 12450
$ws.Range("O4").Value = 1325
$ws.Range("P4").Value = 222
$ws.Range("Q4").Value = 1019
$ws.Range("R4").Value = -1022
$ws.Range("S4").Value = 88
$ws.Range("T4").Value = 692
$ws.Range("U4").Value = 327
$ws.Range("V4").Value = 5023
$ws.Range("W4").Value = 8.960000000000001
$ws.Range("X4").Value = 5.44
$ws.Range("Y4").Value = 9.08
$ws.Range("Z4").Value = 4.73
$ws.Range("AA4").Value = 73.48
$ws.Range("AB4").Value = 5384.08
$ws.Range("AC4").Value = 2482
$ws.Range("AD4").Value = 12.35
$ws.Range("AE4").Value = 28320
$ws.Range("AF4").Value = 1.08
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 0.82
$ws.Range("AI4").Value = 9.99
$ws.Range("AJ4").Value = 44311468

# Row 5
$ws.Range("D5").Value = 20093
$ws.Range("E5").Value = 1811
$ws.Range("F5").Value = 1811
$ws.Range("G5").Value = 1311
$ws.Range("H5").Value = 947
$ws.Range("I5").Value = 1087
$ws.Range("J5").Value = -140
$ws.Range("K5").Value = 22749
$ws.Range("L5").Value = 9341
$ws.Range("M5").Value = 13407
$ws.Range("N5").Value = 12295
$ws.Range("O5").Value = 1112
$ws.Range("P5").Value = 222
$ws.Range("Q5").Value = 1591
$ws.Range("R5").Value = -1673
$ws.Range("S5").Value = -471
$ws.Range("T5").Value = 1443
$ws.Range("U5").Value = 149
$ws.Range("V5").Value = 4329
$ws.Range("W5").Value = 9.01
$ws.Range("X5").Value = 4.71
$ws.Range("Y5").Value = 8.789999999999999
$ws.Range("Z5").Value = 4.06
$ws.Range("AA5").Value = 69.67
$ws.Range("AB5").Value = 5824.27
$ws.Range("AC5").Value = 2453
$ws.Range("AD5").Value = 12.84
$ws.Range("AE5").Value = 28029
$ws.Range("AF5").Value = 1.12
$ws.Range("AG5").Value = 300
$ws.Range("AH5").Value = 0.95
$ws.Range("AI5").Value = 12.11
$ws.Range("AJ5").Value = 44311468

# Row 6
$ws.Range("D6").Value = 21013
$ws.Range("E6").Value = 2010
$ws.Range("F6").Value = 2010
$ws.Range("G6").Value = 1691
$ws.Range("H6").Value = 1113
$ws.Range("I6").Value = 1131
$ws.Range("K6").Value = 26389
$ws.Range("L6").Value = 9758
$ws.Range("M6").Value = 16631
$ws.Range("N6").Value = 15446
$ws.Range("P6").Value = 222
$ws.Range("Q6").Value = 1839
$ws.Range("R6").Value = -1169
$ws.Range("S6").Value = -804
$ws.Range("T6").Value = 1271
$ws.Range("U6").Value = 569
$ws.Range("V6").Value = 3825
$ws.Range("W6").Value = 9.57
$ws.Range("X6").Value = 5.3
$ws.Range("Y6").Value = 8.15
$ws.Range("Z6").Value = 4.53
$ws.Range("AA6").Value = 58.67
$ws.Range("AB6").Value = 6276.46
$ws.Range("AC6").Value = 2551
$ws.Range("AD6").Value = 15.03
$ws.Range("AE6").Value = 35211
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 0.91
$ws.Range("AI6").Value = 13.58
$ws.Range("AJ6").Value = 44311468

# Row 7
$ws.Range("D7").Value = 23647
$ws.Range("E7").Value = 2297
$ws.Range("G7").Value = 2240
$ws.Range("H7").Value = 1687
$ws.Range("I7").Value = 1592
$ws.Range("K7").Value = 30083
$ws.Range("L7").Value = 11131
$ws.Range("M7").Value = 18952
$ws.Range("N7").Value = 17756
$ws.Range("P7").Value = 221
$ws.Range("Q7").Value = 1680
$ws.Range("R7").Value = -531
$ws.Range("S7").Value = 752
$ws.Range("T7").Value = 1128
$ws.Range("U7").Value = 490
$ws.Range("W7").Value = 9.710000000000001
$ws.Range("X7").Value = 7.13
$ws.Range("Y7").Value = 9.59
$ws.Range("Z7").Value = 5.97
$ws.Range("AA7").Value = 58.73
$ws.Range("AC7").Value = 3592
$ws.Range("AD7").Value = 8.6
$ws.Range("AE7").Value = 40479
$ws.Range("AF7").Value = 0.76
$ws.Range("AG7").Value = 361
$ws.Range("AH7").Value = 1.17
$ws.Range("AI7").Value = 10.05

# Row 8
$ws.Range("D8").Value = 24906
$ws.Range("E8").Value = 2616
$ws.Range("G8").Value = 2553
$ws.Range("H8").Value = 1933
$ws.Range("I8").Value = 1841
$ws.Range("K8").Value = 32053
$ws.Range("L8").Value = 11081
$ws.Range("M8").Value = 20970
$ws.Range("N8").Value = 19778
$ws.Range("P8").Value = 221
$ws.Range("Q8").Value = 2380
$ws.Range("R8").Value = -922
$ws.Range("S8").Value = -441
$ws.Range("T8").Value = 1062
$ws.Range("U8").Value = 1109
$ws.Range("W8").Value = 10.5
$ws.Range("X8").Value = 7.76
$ws.Range("Y8").Value = 9.81
$ws.Range("Z8").Value = 6.22
$ws.Range("AA8").Value = 52.84
$ws.Range("AC8").Value = 4155
$ws.Range("AD8").Value = 7.44
$ws.Range("AE8").Value = 45088
$ws.Range("AF8").Value = 0.6899999999999999
$ws.Range("AG8").Value = 372
$ws.Range("AH8").Value = 1.2
$ws.Range("AI8").Value = 8.960000000000001

# Row 9
$ws.Range("D9").Value = 26100
$ws.Range("E9").Value = 2776
$ws.Range("G9").Value = 2758
$ws.Range("H9").Value = 2082
$ws.Range("I9").Value = 2025
$ws.Range("K9").Value = 34306
$ws.Range("L9").Value = 10967
$ws.Range("M9").Value = 23340
$ws.Range("N9").Value = 21974
$ws.Range("P9").Value = 221
$ws.Range("Q9").Value = 2330
$ws.Range("R9").Value = -960
$ws.Range("S9").Value = -399
$ws.Range("T9").Value = 1090
$ws.Range("U9").Value = 1024
$ws.Range("W9").Value = 10.64
$ws.Range("X9").Value = 7.98
$ws.Range("Y9").Value = 9.699999999999999
$ws.Range("Z9").Value = 6.27
$ws.Range("AA9").Value = 46.99
$ws.Range("AC9").Value = 4571
$ws.Range("AD9").Value = 6.76
$ws.Range("AE9").Value = 50093
$ws.Range("AF9").Value = 0.62
$ws.Range("AG9").Value = 362
$ws.Range("AH9").Value = 1.17
$ws.Range("AI9").Value = 7.93
